# Redmine #9229 - Calibration sheet update for GP05MOAS-GL363
# Update anchor launch / recover dates+time, cruise number label, and a
# calibration coefficient value; highlight the edited cells in blue text
# (and the calibration value also gets a yellow fill) to flag the changes.

$wb = $excel.ActiveWorkbook

# --- Sheet "Moorings": update deployment row (row 2) ---
$ws1 = $wb.Worksheets.Item("Moorings")

# Anchor Launch Date: 1-Jul-2013 -> 19-Jul-2013
$ws1.Range("D2").Value = [DateTime]::new(2013, 7, 19)
# Anchor Launch Time: 00:00 -> 23:00
$ws1.Range("E2").Value = [DateTime]::new(1899, 12, 30, 23, 0, 0)
# Recover Date: 13-Jun-2014 -> 15-Jun-2014
$ws1.Range("F2").Value = [DateTime]::new(2014, 6, 15)
# Cruise Number: "Melville 130" -> "MV-130"
$ws1.Range("J2").Value = "MV-130"

# Flag the edited cells with blue font color (RGB 0,0,255)
$ws1.Range("D2").Font.Color = 255 * 65536 * 0 + 255 * 256 * 0 + 255
$ws1.Range("D2,E2,F2,J2").Font.Color = 16711680 -as [int]
$ws1.Range("D2,E2,F2,J2").Font.Color = 0x0000FF

# --- Sheet "Asset_Cal_Info": update CC_angular_resolution value (row 6) ---
$ws2 = $wb.Worksheets.Item("Asset_Cal_Info")
$ws2.Range("F6").Value = 1.096
$ws2.Range("F6").Font.Color = 0x0000FF
$ws2.Range("F6").Interior.Color = 0x00FFFF
